# Natmi following Dr Hou advice
# Rebuild the Sending/Target cluster cross-product (ECs x FAPs x sCs) for the Hgf-Sdc1 pair.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs (Hgf/Sdc1)
$row = New-Object 'object[,]' 1,20
$row[0,0] = "ECs"
$row[0,1] = "Hgf"
$row[0,2] = "Sdc1"
$row[0,3] = "ECs"
$row[0,4] = 2
$row[0,5] = 0.6666666666666666
$row[0,6] = 5.928568666666667
$row[0,7] = 17.785706
$row[0,8] = 0.3809768389628236
$row[0,9] = 0.3809768389628236
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 0.8369776666666665
$row[0,13] = 2.510933
$row[0,14] = 0.0694586718035551
$row[0,15] = 0.06945867180355511
$row[0,16] = 4.962079569299777
$row[0,17] = 44.658716123698
$row[0,18] = 0.02646214522227463
$row[0,19] = 0.02646214522227464
$ws.Range("A2:T2").Value = $row

# Row 3: ECs -> FAPs (Hgf/Sdc1)
$row = New-Object 'object[,]' 1,20
$row[0,0] = "ECs"
$row[0,1] = "Hgf"
$row[0,2] = "Sdc1"
$row[0,3] = "FAPs"
$row[0,4] = 2
$row[0,5] = 0.6666666666666666
$row[0,6] = 5.928568666666667
$row[0,7] = 17.785706
$row[0,8] = 0.3809768389628236
$row[0,9] = 0.3809768389628236
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 2.427350333333333
$row[0,13] = 7.282051
$row[0,14] = 0.2014397000898671
$row[0,15] = 0.2014397000898671
$row[0,16] = 14.39071312922289
$row[0,17] = 129.516418163006
$row[0,18] = 0.0767438601818568
$row[0,19] = 0.0767438601818568
$ws.Range("A3:T3").Value = $row

# Row 4: ECs -> sCs (Hgf/Sdc1)
$row = New-Object 'object[,]' 1,20
$row[0,0] = "ECs"
$row[0,1] = "Hgf"
$row[0,2] = "Sdc1"
$row[0,3] = "sCs"
$row[0,4] = 2
$row[0,5] = 0.6666666666666666
$row[0,6] = 5.928568666666667
$row[0,7] = 17.785706
$row[0,8] = 0.3809768389628236
$row[0,9] = 0.3809768389628236
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 8.785681666666667
$row[0,13] = 26.357045
$row[0,14] = 0.7291016281065776
$row[0,15] = 0.7291016281065776
$row[0,16] = 52.08651704430778
$row[0,17] = 468.77865339877
$row[0,18] = 0.2777708335586921
$row[0,19] = 0.2777708335586921
$ws.Range("A4:T4").Value = $row

# Row 5: FAPs -> ECs (Hgf/Sdc1)
$row = New-Object 'object[,]' 1,20
$row[0,0] = "FAPs"
$row[0,1] = "Hgf"
$row[0,2] = "Sdc1"
$row[0,3] = "ECs"
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 8.583521
$row[0,7] = 25.750563
$row[0,8] = 0.5515872180307627
$row[0,9] = 0.5515872180307626
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 0.8369776666666665
$row[0,13] = 2.510933
$row[0,14] = 0.0694586718035551
$row[0,15] = 0.06945867180355511
$row[0,16] = 7.184215378364331
$row[0,17] = 64.65793840527898
$row[0,18] = 0.03831251554823473
$row[0,19] = 0.03831251554823473
$ws.Range("A5:T5").Value = $row

# Row 6: FAPs -> FAPs (Hgf/Sdc1)
$row = New-Object 'object[,]' 1,20
$row[0,0] = "FAPs"
$row[0,1] = "Hgf"
$row[0,2] = "Sdc1"
$row[0,3] = "FAPs"
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 8.583521
$row[0,7] = 25.750563
$row[0,8] = 0.5515872180307627
$row[0,9] = 0.5515872180307626
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 2.427350333333333
$row[0,13] = 7.282051
$row[0,14] = 0.2014397000898671
$row[0,15] = 0.2014397000898671
$row[0,16] = 20.83521256052367
$row[0,17] = 187.516913044713
$row[0,18] = 0.111111563773521
$row[0,19] = 0.111111563773521
$ws.Range("A6:T6").Value = $row

# Row 7: FAPs -> sCs (Hgf/Sdc1)
$row = New-Object 'object[,]' 1,20
$row[0,0] = "FAPs"
$row[0,1] = "Hgf"
$row[0,2] = "Sdc1"
$row[0,3] = "sCs"
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 8.583521
$row[0,7] = 25.750563
$row[0,8] = 0.5515872180307627
$row[0,9] = 0.5515872180307626
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 8.785681666666667
$row[0,13] = 26.357045
$row[0,14] = 0.7291016281065776
$row[0,15] = 0.7291016281065776
$row[0,16] = 75.41208308514832
$row[0,17] = 678.708747766335
$row[0,18] = 0.4021631387090069
$row[0,19] = 0.4021631387090068
$ws.Range("A7:T7").Value = $row

# Row 8: sCs -> ECs (Hgf/Sdc1)
$row = New-Object 'object[,]' 1,20
$row[0,0] = "sCs"
$row[0,1] = "Hgf"
$row[0,2] = "Sdc1"
$row[0,3] = "ECs"
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 1.049404
$row[0,7] = 3.148212
$row[0,8] = 0.06743594300641363
$row[0,9] = 0.06743594300641362
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 0.8369776666666665
$row[0,13] = 2.510933
$row[0,14] = 0.0694586718035551
$row[0,15] = 0.06945867180355511
$row[0,16] = 0.8783277113106666
$row[0,17] = 7.904949401795998
$row[0,18] = 0.004684011033045731
$row[0,19] = 0.004684011033045731
$ws.Range("A8:T8").Value = $row

# Row 9: sCs -> FAPs (Hgf/Sdc1)
$row = New-Object 'object[,]' 1,20
$row[0,0] = "sCs"
$row[0,1] = "Hgf"
$row[0,2] = "Sdc1"
$row[0,3] = "FAPs"
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 1.049404
$row[0,7] = 3.148212
$row[0,8] = 0.06743594300641363
$row[0,9] = 0.06743594300641362
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 2.427350333333333
$row[0,13] = 7.282051
$row[0,14] = 0.2014397000898671
$row[0,15] = 0.2014397000898671
$row[0,16] = 2.547271149201333
$row[0,17] = 22.925440342812
$row[0,18] = 0.01358427613448934
$row[0,19] = 0.01358427613448933
$ws.Range("A9:T9").Value = $row

# Row 10: sCs -> sCs (Hgf/Sdc1)
$row = New-Object 'object[,]' 1,20
$row[0,0] = "sCs"
$row[0,1] = "Hgf"
$row[0,2] = "Sdc1"
$row[0,3] = "sCs"
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 1.049404
$row[0,7] = 3.148212
$row[0,8] = 0.06743594300641363
$row[0,9] = 0.06743594300641362
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 8.785681666666667
$row[0,13] = 26.357045
$row[0,14] = 0.7291016281065776
$row[0,15] = 0.7291016281065776
$row[0,16] = 9.219729483726667
$row[0,17] = 82.97756535354
$row[0,18] = 0.04916765583887855
$row[0,19] = 0.04916765583887855
$ws.Range("A10:T10").Value = $row
